$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.933.41'
$ws.Range("E2").Value = '  +0.04%  '
$ws.Range("D3").Value = '1.630.03'
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").Value = "'211.70"
$ws.Range("E5").Value = '  -0.91%  '
$ws.Range("E6").Value = '  -0.10%  '
$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("D8").Value = "'23.38"
$ws.Range("E8").Value = '  -1.75%  '
$ws.Range("E9").Value = '  -2.11%  '
$ws.Range("E11").Value = '  +0.50%  '
$ws.Range("D12").Value = '1.860.40'
$ws.Range("E12").Value = '  -0.74%  '
$ws.Range("D13").Value = '1.622.00'
$ws.Range("E13").Value = '  -1.21%  '
$ws.Range("E14").Value = '  -1.20%  '
$ws.Range("D15").Value = "'0.562"
$ws.Range("E15").Value = '  -3.11%  '
$ws.Range("D16").Value = "'65.58"
$ws.Range("E16").Value = '  -0.89%  '
$ws.Range("D17").Value = '27.929.95'
$ws.Range("E17").Value = '  +0.06%  '
$ws.Range("D18").Value = "'230.63"
$ws.Range("E18").Value = '  -0.45%  '
$ws.Range("E19").Value = '  -0.20%  '
$ws.Range("E20").Value = '  +0.04%  '
$ws.Range("E21").Value = '  -0.10%  '
$ws.Range("E22").Value = '  -9.97%  '
$ws.Range("E23").Value = '  -0.98%  '
$ws.Range("D24").Value = "'2.04"
$ws.Range("E24").Value = '  -1.30%  '
$ws.Range("E25").Value = '  +2.09%  '
$ws.Range("D26").Value = "'6.93"
$ws.Range("E26").Value = '  -0.33%  '
$ws.Range("E27").Value = '  -0.61%  '
$ws.Range("D28").Value = "'15.57"
$ws.Range("E28").Value = '  -1.14%  '
$ws.Range("E29").Value = '  -0.03%  '
$ws.Range("E30").Value = '  -0.74%  '
$ws.Range("D31").Value = "'0.0482"
$ws.Range("E31").Value = '  -0.89%  '
$ws.Range("E32").Value = '  +1.55%  '
$ws.Range("E33").Value = '  -1.73%  '
$ws.Range("D34").Value = '1.396.45'
$ws.Range("E34").Value = '  -2.18%  '
$ws.Range("E35").Value = '  -0.30%  '
$ws.Range("D36").Value = "'1.04"
$ws.Range("E36").Value = '  +13.01%  '
$ws.Range("E37").Value = '  +0.54%  '
$ws.Range("E38").Value = '  +1.94%  '
$ws.Range("E39").Value = '  -0.39%  '
$ws.Range("D40").Value = "'0.865"
$ws.Range("E40").Value = '  -3.17%  '
$ws.Range("E41").Value = '  -0.61%  '
$ws.Range("E42").Value = '  -0.05%  '
$ws.Range("D43").Value = "'66.48"
$ws.Range("E43").Value = '  -1.32%  '
$ws.Range("E44").Value = '  -0.32%  '
$ws.Range("D45").Value = "'5.45"
$ws.Range("E45").Value = '  -0.11%  '
$ws.Range("E46").Value = '  -0.53%  '
$ws.Range("D47").Value = '1.770.95'
$ws.Range("E47").Value = '  -0.67%  '
$ws.Range("D48").Value = "'88.11"
$ws.Range("E48").Value = '  -0.58%  '
$ws.Range("D49").Value = '0.0₆0104'
$ws.Range("E49").Value = '  -0.59%  '
$ws.Range("E50").Value = '  -1.23%  '
$ws.Range("E51").Value = '  -0.43%  '
